# Updates the "广州-漫展信息" workbook to the state generated at commit 456a3b4.
# Sheet order: 1 = 展览 (exhibitions), 2 = 演出 (performances),
#              3 = 本地生活 (local life), 4 = 全部类型 (all types)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Simple "想去人数" (want-to-go count) bumps that don't move any rows.
$ws1.Range("F2").Value = 2494
$ws1.Range("F4").Value = 263
$ws1.Range("F5").Value = 429
$ws1.Range("F6").Value = 740
$ws1.Range("F8").Value = 932
$ws1.Range("F10").Value = 988
$ws1.Range("F16").Value = 1129

# Row 17: count bump and the ticket status flips from the "已售罄" (sold out)
# text label to a plain 0 (numeric) in the lowest-price column.
$ws1.Range("F17").Value = 24944
$ws1.Range("G17").Value = 0

$ws1.Range("F18").Value = 2450
$ws1.Range("F20").Value = 373
$ws1.Range("F22").Value = 106
$ws1.Range("F25").Value = 135
$ws1.Range("F28").Value = 90
$ws1.Range("F30").Value = 381

# A brand-new exhibition is inserted as row 31, pushing the former rows
# 31-33 down to 32-34.
$ws1.Rows.Item(31).Insert()

# The inserted row's index cell (column A) needs the same look as the other
# index cells - easiest way is to copy the formatting from a neighbour and
# then overwrite the value.
$ws1.Range("A32").Copy($ws1.Range("A31"))

$ws1.Range("A31").Value = 30
$ws1.Range("B31").NumberFormat = "@"
$ws1.Range("B31").Value = "2024-08-17"
$ws1.Range("B31").ClearFormats()
$ws1.Range("C31").Value = "广州·COC星火次元云漫创作交流展"
$ws1.Range("D31").Value = "黄边三横路一街1号 设计殿堂"
$ws1.Range("E31").Value = "2024.08.17 09:00-08.18 17:00"
$ws1.Range("F31").Value = 0
$ws1.Range("G31").Value = 68
$ws1.Range("H31").Value = "https://show.bilibili.com/platform/detail.html?id=87777"
$ws1.Range("I31").Value = "//i1.hdslb.com/bfs/openplatform/202406/nVPxhUFQ1718936306088.jpeg"

# The index column (A) holds plain numbers, not formulas, so the rows that
# got shifted down by the insert need their index re-stamped to row-1.
$ws1.Range("A32").Value = 31
$ws1.Range("A33").Value = 32
$ws1.Range("A34").Value = 33

# ---------------------------------------------------------------------------
# Sheet 2: 演出
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("F6").Value = 285
$ws2.Range("F7").Value = 248
$ws2.Range("F8").Value = 33

# Row 10: same sold-out-label -> numeric-price flip as sheet 1 row 17.
$ws2.Range("F10").Value = 3699
$ws2.Range("G10").Value = 480

$ws2.Range("F18").Value = 30
$ws2.Range("F21").Value = 4148

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("F2").Value = 283
$ws3.Range("F3").Value = 181
$ws3.Range("F4").Value = 843

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2").Value = 283
$ws4.Range("F3").Value = 181
$ws4.Range("F4").Value = 2494
$ws4.Range("F5").Value = 843
$ws4.Range("F7").Value = 263
$ws4.Range("F8").Value = 429
$ws4.Range("F9").Value = 740
$ws4.Range("F14").Value = 285
$ws4.Range("F16").Value = 932
$ws4.Range("F18").Value = 988
$ws4.Range("F23").Value = 1129

# Row 24: count bump; status text changes from "已售罄" to "暂时售罄" but
# stays a text label (unlike the analogous rows on sheets 1 and 2).
$ws4.Range("F24").Value = 24944
$ws4.Range("G24").Value = "暂时售罄"

$ws4.Range("F25").Value = 33
$ws4.Range("F31").Value = 2451
$ws4.Range("F34").Value = 373
$ws4.Range("F41").Value = 30
$ws4.Range("F42").Value = 90
